$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.485.02"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "1.989.77"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4677"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3947"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07939"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("E12").Value = "  +5.07%  "
$ws.Range("D13").Value = "2.008.92"
$ws.Range("E13").Value = "  +7.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.266"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.877"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009961"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "29.555.33"
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.545"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.097"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.992"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "120.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.953"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09442"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9133"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.351"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.254"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.000003511"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +102.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05843"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.173"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.890"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5773"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1829"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.829"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5381"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.166"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.872"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06939"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3079"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.28%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.779"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.40%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.85%  "
